# Applies the "additional formatting changes" commit:
#   - Heading4 style: paragraph spacing before/after 200/0 -> 120/120 (6pt/6pt)
#   - Heading7 style: paragraph spacing before/after 200/0 -> 320/120 (16pt/6pt),
#                      plus 1.5 line spacing (line=360, auto), and bold added
#   - Heading4Char / Heading7Char (linked character styles) mirror the font change
#   - Caption style: add explicit spacing before/after 240/240 (12pt/12pt)
#   - CaptionChar (linked character style) stays in sync

$d = $word.ActiveDocument

# --- Heading 4 -----------------------------------------------------------
$heading4 = $d.Styles.Item("Heading4")
$heading4.ParagraphFormat.SpaceBefore = 6
$heading4.ParagraphFormat.SpaceAfter = 6

# --- Heading 7 -------------------------------------------------------------
$heading7 = $d.Styles.Item("Heading7")
$heading7.ParagraphFormat.SpaceBefore = 16
$heading7.ParagraphFormat.SpaceAfter = 6
$heading7.ParagraphFormat.LineSpacingRule = 5   # wdLineSpaceMultiple
$heading7.ParagraphFormat.LineSpacing = 18      # 360 twips / 20 -> 1.5 lines at 12pt base
$heading7.Font.Bold = $true

# Linked character style ("Heading 7 Char") must pick up the new bold weight too.
$heading7Char = $d.Styles.Item("Heading7Char")
$heading7Char.Font.Bold = $true

# --- Caption -----------------------------------------------------------
$caption = $d.Styles.Item("Caption")
$caption.ParagraphFormat.SpaceBefore = 12
$caption.ParagraphFormat.SpaceAfter = 12

Write-Output "Styles updated: Heading4, Heading7, Heading7Char, Caption"
